$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2541.2
$ws.Range("I28").Value = 1065
$ws.Range("K28").Value = 1065
$ws.Range("M28").Value = -580

$ws.Range("H33").Value = 58824490
$ws.Range("I33").Value = 459.18182
$ws.Range("J33").Value = 166668540
$ws.Range("K33").Value = 459.18182
$ws.Range("L33").Value = 166668540
$ws.Range("M33").Value = -230.18182
$ws.Range("N33").Value = -166668998

$ws.Range("H43").Value = 886
$ws.Range("I43").Value = 714.2857
$ws.Range("J43").Value = 1057.7142
$ws.Range("K43").Value = 714.2857
$ws.Range("L43").Value = 1057.7142
$ws.Range("M43").Value = -645.2857
$ws.Range("N43").Value = -1195.7142

$ws.Range("H129").Value = 1054.5
$ws.Range("J129").Value = 1179.125
$ws.Range("L129").Value = 3537.375
$ws.Range("N129").Value = -13537.375

$ws.Range("H132").Value = 3238.1
$ws.Range("I132").Value = 1274.9678
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 3824.9034
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -1294.9034
$ws.Range("N132").Value = -35060

$ws.Range("H138").Value = 2229.83
$ws.Range("I138").Value = 917.2
$ws.Range("J138").Value = 3104.9167
$ws.Range("K138").Value = 2751.6
$ws.Range("L138").Value = 9314.750100000001
$ws.Range("M138").Value = 2388.4
$ws.Range("N138").Value = -19594.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10227
$ws.Range("J43").Value = 10227
$ws.Range("L43").Value = 10227
$ws.Range("N43").Value = -10853

$ws.Range("H63").Value = 3580
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 4300
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 4300
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -5672

$ws.Range("H66").Value = 3580
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 4300
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 21500
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -28364

$ws.Range("H74").Value = 29662.4
$ws.Range("I74").Value = 35468.1
$ws.Range("J74").Value = 1601.5
$ws.Range("K74").Value = 35468.1
$ws.Range("L74").Value = 1601.5
$ws.Range("M74").Value = -34594.1
$ws.Range("N74").Value = -3349.5

$ws.Range("H77").Value = 29662.4
$ws.Range("I77").Value = 35468.1
$ws.Range("J77").Value = 1601.5
$ws.Range("K77").Value = 177340.5
$ws.Range("L77").Value = 8007.5
$ws.Range("M77").Value = -172972.5
$ws.Range("N77").Value = -16743.5

$ws.Range("H97").Value = 1250
$ws.Range("I97").Value = 1190.909
$ws.Range("J97").Value = 1412.5
$ws.Range("K97").Value = 1190.909
$ws.Range("L97").Value = 1412.5
$ws.Range("M97").Value = -694.9090000000001
$ws.Range("N97").Value = -2404.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1354.4615
$ws.Range("I94").Value = 1289.8
$ws.Range("J94").Value = 1570
$ws.Range("K94").Value = 1289.8
$ws.Range("L94").Value = 1570
$ws.Range("M94").Value = -838.8
$ws.Range("N94").Value = -2472

$ws.Range("H134").Value = 5587.85
$ws.Range("I134").Value = 5767.875
$ws.Range("J134").Value = 5467.8335
$ws.Range("K134").Value = 17303.625
$ws.Range("L134").Value = 16403.5005
$ws.Range("M134").Value = -14768.625
$ws.Range("N134").Value = -21473.5005

$ws.Range("H141").Value = 85695
$ws.Range("J141").Value = 85695
$ws.Range("L141").Value = 85695
$ws.Range("N141").Value = -96055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 630
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 1500
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 4500
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -4780

$ws.Range("H17").Value = 5750
$ws.Range("I17").Value = 999.3333
$ws.Range("J17").Value = 20002
$ws.Range("K17").Value = 2997.9999
$ws.Range("L17").Value = 60006
$ws.Range("M17").Value = -2828.9999
$ws.Range("N17").Value = -60344

$ws.Range("H98").Value = 1489.7
$ws.Range("I98").Value = 1987.25
$ws.Range("J98").Value = 1158
$ws.Range("K98").Value = 5961.75
$ws.Range("L98").Value = 3474
$ws.Range("M98").Value = -4463.75
$ws.Range("N98").Value = -6470

$ws.Range("H100").Value = 4385.7144
$ws.Range("J100").Value = 4385.7144
$ws.Range("L100").Value = 13157.1432
$ws.Range("N100").Value = -14779.1432

$ws.Range("H109").Value = 2463.4
$ws.Range("I109").Value = 1613.5
$ws.Range("J109").Value = 3030
$ws.Range("K109").Value = 4840.5
$ws.Range("L109").Value = 9090
$ws.Range("M109").Value = -3800.5
$ws.Range("N109").Value = -11170

$ws.Range("H131").Value = 848.47
$ws.Range("J131").Value = 883.9659
$ws.Range("L131").Value = 2651.8977
$ws.Range("N131").Value = -12731.8977

$ws.Range("H132").Value = 3487918.8
$ws.Range("I132").Value = 1627928.2
$ws.Range("J132").Value = 11113880
$ws.Range("K132").Value = 14651353.8
$ws.Range("L132").Value = 100024920
$ws.Range("M132").Value = -14648823.8
$ws.Range("N132").Value = -100029980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1382.6316
$ws.Range("I97").Value = 1007
$ws.Range("K97").Value = 1007
$ws.Range("M97").Value = -511

$ws.Range("H102").Value = 1344.5416
$ws.Range("I102").Value = 1102.6875
$ws.Range("J102").Value = 1828.25
$ws.Range("K102").Value = 1102.6875
$ws.Range("L102").Value = 1828.25
$ws.Range("M102").Value = 519.3125
$ws.Range("N102").Value = -5072.25

$ws.Range("H122").Value = 47780.41
$ws.Range("I122").Value = 64404.375
$ws.Range("J122").Value = 3449.8333
$ws.Range("K122").Value = 193213.125
$ws.Range("L122").Value = 10349.4999
$ws.Range("M122").Value = -190763.125
$ws.Range("N122").Value = -15249.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1237.4546
$ws.Range("I82").Value = 1714.5714
$ws.Range("J82").Value = 1014.8
$ws.Range("K82").Value = 1714.5714
$ws.Range("L82").Value = 1014.8
$ws.Range("M82").Value = -1353.5714
$ws.Range("N82").Value = -1736.8

$ws.Range("H85").Value = 1237.4546
$ws.Range("I85").Value = 1714.5714
$ws.Range("J85").Value = 1014.8
$ws.Range("K85").Value = 1714.5714
$ws.Range("L85").Value = 1014.8
$ws.Range("M85").Value = -466.5714
$ws.Range("N85").Value = -3510.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -627
$ws.Range("N96").ClearContents()

$ws.Range("H126").Value = 83334050
$ws.Range("I126").Value = 90909780
$ws.Range("K126").Value = 272729340
$ws.Range("M126").Value = -272726870
